# Flip the Runmode for the "ValidateCRMTest" test case from "N" to "Y"
# on both the TestCases sheet (the master switch) and the TestData sheet
# (the per-test run flag).

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestCases.Range("B3").Value = "Y"

$wsTestData = $wb.Worksheets.Item("TestData")
$wsTestData.Range("A9").Value = "Y"

# Leave the cursor sitting just past the last used row, mirroring the
# author's final selection after editing the table.
$wsTestData.Activate()
$wsTestData.Range("B10").Select()
